# Apply the diff: add two new rows (f_1_t / n_1_t) referencing existing
# descriptions, with a new highlighted style (Palatino 16pt, custom color)
# on the first of the two new rows, and update the view selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 197: "f_1_t" / "feminine noun/adjective ending in t" -- highlighted
$ws.Range("A197").Value = "f_1_t"
$ws.Range("B197").Value = "feminine noun/adjective ending in t"

$highlightRange = $ws.Range("A197:B197")
$highlightRange.Font.Name = "Palatino"
$highlightRange.Font.Size = 16
$highlightRange.Font.Color = 5722185
$ws.Rows.Item(197).RowHeight = 21

# New row 198: "n_1_t" / "neuter noun/adjective ending in t"
$ws.Range("A198").Value = "n_1_t"
$ws.Range("B198").Value = "neuter noun/adjective ending in t"

# Update the view: scroll position + active selection
$ws.Application.ActiveWindow.ScrollRow = 161
$ws.Range("E182").Select()
